$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: full path -> just the file name
$ws.Range("B1").Value = "SimpleSampleCode"

# B2: "put number here" -> number 3
$ws.Range("B2").Value = 3

# B3: 3.0 -> number 0
$ws.Range("B3").Value = 0

# B4: "put number here" -> number 0
$ws.Range("B4").Value = 0
